{"js": "// Fix the capitalization of \"Ejercicio Gen\u00e9rico\" -> \"Ejercicio gen\u00e9rico\" in the\n// document title and move the \"_GoBack\" bookmark (Word's \"last edit position\"\n// marker) from its old spot (right after \"CS_04_02_CO \") to the point in the\n// title where the edit was made, splitting the original single run into two\n// runs around the bookmark - exactly mirroring what Word itself does when a\n// user edits text and re-saves.\n\nconst body = context.document.body;\n\n// 1) Remove the old \"_GoBack\" bookmark (harmless no-op if it is not present).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Locate \"Ejercicio G\" (case sensitive) - the first 11 characters of the\n//    title run, i.e. everything up to and including the capital \"G\" that\n//    needs to become lower-case.\nconst found = body.search(\"Ejercicio G\", { matchCase: true, matchWholeWord: false });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length > 0) {\n  const titleStart = found.items[0];\n\n  // 3) Drop a bookmark right after \"Ejercicio G\" - this is the collapsed\n  //    point where the rest of the run (\"en\u00e9rico \") will keep living; Word\n  //    always splits the underlying run around a newly-inserted bookmark,\n  //    which is exactly the run layout the target document has.\n  const boundary = titleStart.getRange(\"End\");\n  boundary.insertBookmark(\"_GoBack\");\n  await context.sync();\n\n  // 4) Fix the capitalization: \"Ejercicio G\" -> \"Ejercicio g\".\n  titleStart.insertText(\"Ejercicio g\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Fix the capitalization of \"Ejercicio Gen\u00e9rico\" -> \"Ejercicio gen\u00e9rico\" in the\n# document title, and relocate the \"_GoBack\" bookmark (Word's \"last edit\n# position\" marker) from its old spot (right after \"CS_04_02_CO \") to the\n# point in the title where the edit was made. Adding a bookmark with a name\n# that is already in use moves it (Word enforces unique bookmark names), so\n# re-adding \"_GoBack\" at the new location automatically removes it from the\n# old one - exactly mirroring what a real editing session does.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.Text = \"Ejercicio G\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$found = $rng.Find.Execute()\n\nif ($found) {\n    # Drop the \"_GoBack\" bookmark right after \"Ejercicio G\" (collapsed point);\n    # Word splits the underlying run around a newly-inserted bookmark, giving\n    # the exact two-run layout (\"Ejercicio g\" | bookmark | \"en\u00e9rico \") the\n    # target document has.\n    $bmRng = $rng.Duplicate\n    $bmRng.Collapse(0)   # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $bmRng)\n\n    # Fix the capitalization: \"Ejercicio G\" -> \"Ejercicio g\".\n    $rng.Text = \"Ejercicio g\"\n}\n"}
